$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C57 previously held "NA"; the refreshed script run found no page number,
# so the cell is now blank.
$ws.Range("C57").NumberFormat = "@"
$ws.Range("C57").Value = ""
$ws.Range("C57").Style = "Normal"

# Append the newest scraped row (row 58) with the new search-term result.
$ws.Range("A58").NumberFormat = "@"
$ws.Range("A58").Value = "2025-04-24"
$ws.Range("A58").Style = "Normal"

$ws.Range("B58").Value = "buse"
$ws.Range("C58").Value = 67
$ws.Range("D58").Value = 1
